# "Generate Report for Handback"
#
# The localization-status report is re-generated after a handback: the
# Overview "Status" column moves from "Ready for handoff" to "Handed back:
# in sync with en-US", and each locale sheet (zh-cn, de-de) gets its
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns populated (with a hyperlink on the target-file cell,
# matching the look of the existing source-file hyperlink in column A).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$hyperlinkColor = 15570276   # RGB(0x64,0x95,0xED) == the workbook's existing "FF6495ED" hyperlink color

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9024f037eb97c717e47d43887cb9b563224388f/e2e/"

$file1 = "8bced00c-8463-4e94-aabf-6df5b9c8127c"
$file2 = "dca63bb7-9d6e-4b08-b3cf-5c7f206b255e"

# ---------------------------------------------------------------------
# 1) Overview sheet: refresh the per-locale status cells
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------
# 2) zh-cn sheet: Status + Latest Target File / Handback File / Handback DateTime
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("I2").Value = ($file1 + ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), ($repoBase + $file1 + ".md"), "", "", ($file1 + ".md"))
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = $hyperlinkColor
$zhcn.Range("J2").Value = ($file1 + ".a667a8ab4a53204590636f03be9eabbfc8ecff33.zh-cn.xlf")
$zhcn.Range("K2").Value = "2016-08-15 12:24:49"

$zhcn.Range("I3").Value = ($file2 + ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), ($repoBase + $file2 + ".md"), "", "", ($file2 + ".md"))
$zhcn.Range("I3").Font.Underline = 2
$zhcn.Range("I3").Font.Color = $hyperlinkColor
$zhcn.Range("J3").Value = ($file2 + ".0cb3d524e4a8d0fdf927f393d6368031ed6fc7e6.zh-cn.xlf")
$zhcn.Range("K3").Value = "2016-08-15 12:24:49"

$zhcn.Columns.Item(3).ColumnWidth = 29.17
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# 3) de-de sheet: Status + Latest Target File / Handback File / Handback DateTime
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("I2").Value = ($file1 + ".md")
$dede.Hyperlinks.Add($dede.Range("I2"), ($repoBase + $file1 + ".md"), "", "", ($file1 + ".md"))
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = $hyperlinkColor
$dede.Range("J2").Value = ($file1 + ".a667a8ab4a53204590636f03be9eabbfc8ecff33.de-de.xlf")
$dede.Range("K2").Value = "2016-08-15 12:24:57"

$dede.Range("I3").Value = ($file2 + ".md")
$dede.Hyperlinks.Add($dede.Range("I3"), ($repoBase + $file2 + ".md"), "", "", ($file2 + ".md"))
$dede.Range("I3").Font.Underline = 2
$dede.Range("I3").Font.Color = $hyperlinkColor
$dede.Range("J3").Value = ($file2 + ".0cb3d524e4a8d0fdf927f393d6368031ed6fc7e6.de-de.xlf")
$dede.Range("K3").Value = "2016-08-15 12:24:57"

$dede.Columns.Item(3).ColumnWidth = 29.17
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17
